# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# New values for column F, keyed by row number. Identical across both sheets.
$updates = @{
    3  = 507
    5  = 8509
    6  = 335
    7  = 1505
    9  = 122
    11 = 250
    12 = 385
    13 = 242
    18 = 453
    19 = 1226
    20 = 184
    21 = 77
    22 = 133
    23 = 93
    25 = 69
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
